$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    85  = "first day of the month"
    86  = "second day of the month"
    87  = "third day of the month"
    88  = "fourth day of the month"
    89  = "fifth day of the month"
    90  = "sixth day of the month"
    91  = "seventh day of the month"
    92  = "eighth day of the month"
    93  = "ninth day of the month"
    94  = "tenth day of the month"
    95  = "eleventh day of the month"
    96  = "twelth day of the month"
    97  = "thirteenth day of the month"
    98  = "fourteenth day of the month"
    99  = "fifteenth day of the month"
    100 = "sixteenth day of the month"
    101 = "seventeenth day of the month"
    102 = "eighteenth day of the month"
    103 = "nineteenth day of the month"
    104 = "twentieth day of the month"
    105 = "twenty-first day of the month"
    106 = "twenty-second day of the month"
    107 = "twenty-third day of the month"
    108 = "twenty-fourth day of the month"
    109 = "twenty-fifth day of the month"
    110 = "twenty-sixth day of the month"
    111 = "twenty-seventh day of the month"
    112 = "twenty-eighth day of the month"
    113 = "twenty-ninth day of the month"
    114 = "thirtieth day of the month"
    115 = "thiry-first day of the month"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
